$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9 (G=5487)
$ws.Range("H9").Value = 243.71428
$ws.Range("I9").Value = 250.83333
$ws.Range("J9").Value = 201
$ws.Range("K9").Value = 250.83333
$ws.Range("L9").Value = 201
$ws.Range("M9").Value = -81.83332999999999
$ws.Range("N9").Value = -539

# Row 40 (G=5505)
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 5000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -5350

# Row 88 (G=12608)
$ws.Range("H88").Value = 1808.2222
$ws.Range("I88").Value = 1775.6
$ws.Range("K88").Value = 1775.6
$ws.Range("M88").Value = -1369.6

# Row 91 (G=12608)
$ws.Range("H91").Value = 1808.2222
$ws.Range("I91").Value = 1775.6
$ws.Range("K91").Value = 1775.6
$ws.Range("M91").Value = -371.5999999999999

# Row 112 (G=27960)
$ws.Range("H112").Value = 1612
$ws.Range("J112").Value = 2294.5
$ws.Range("L112").Value = 6883.5
$ws.Range("N112").Value = -9099.5

# Row 127 (G=36114)
$ws.Range("H127").Value = 362189.8
$ws.Range("I127").Value = 450737.5
$ws.Range("K127").Value = 1352212.5
$ws.Range("M127").Value = -1347252.5

# Row 131 (G=36108)
$ws.Range("H131").Value = 1199.9117
$ws.Range("I131").Value = 1199.9117
$ws.Range("K131").Value = 3599.7351
$ws.Range("M131").Value = 1440.2649

# Row 134 (G=41997)
$ws.Range("H134").Value = 107225.73
$ws.Range("J134").Value = 99998.55499999999
$ws.Range("L134").Value = 99998.55499999999
$ws.Range("N134").Value = -110138.555

# Row 137 (G=44013)
$ws.Range("H137").Value = 3536
$ws.Range("I137").Value = 3451.4666
$ws.Range("J137").Value = 3641.6667
$ws.Range("K137").Value = 10354.3998
$ws.Range("L137").Value = 10925.0001
$ws.Range("M137").Value = -7804.399800000001
$ws.Range("N137").Value = -16025.0001

# Row 141 (G=44161)
$ws.Range("H141").Value = 22938.303
$ws.Range("I141").Value = 22938.303
$ws.Range("K141").Value = 68814.909
$ws.Range("M141").Value = -63634.909

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (G=27713)
$ws.Range("H2").Value = 2532.3157
$ws.Range("I2").Value = 2579.7334
$ws.Range("J2").Value = 2354.5
$ws.Range("K2").Value = 2579.7334
$ws.Range("L2").Value = 2354.5
$ws.Range("M2").Value = -2466.7334
$ws.Range("N2").Value = -2580.5

# Row 32 (G=44147)
$ws.Range("H32").Value = 18006.836
$ws.Range("I32").Value = 2670.5967
$ws.Range("J32").Value = 104447.45
$ws.Range("K32").Value = 2670.5967
$ws.Range("L32").Value = 104447.45
$ws.Range("M32").Value = -2383.5967
$ws.Range("N32").Value = -105021.45

# Row 45 (G=27714)
$ws.Range("H45").Value = 441954.53
$ws.Range("I45").Value = 507482.56
$ws.Range("K45").Value = 507482.56
$ws.Range("M45").Value = -507105.56

# Row 61 (G=43999)
$ws.Range("H61").Value = 2427
$ws.Range("I61").Value = 2427
$ws.Range("K61").Value = 2427
$ws.Range("M61").Value = -2215

# Row 74 (G=44000)
$ws.Range("H74").Value = 2501.8333
$ws.Range("I74").Value = 2547.4546
$ws.Range("K74").Value = 2547.4546
$ws.Range("M74").Value = -1673.4546

# Row 77 (G=44000)
$ws.Range("H77").Value = 2501.8333
$ws.Range("I77").Value = 2547.4546
$ws.Range("K77").Value = 12737.273
$ws.Range("M77").Value = -8369.273000000001

# Row 110 (G=27708)
$ws.Range("H110").Value = 6876.375
$ws.Range("I110").Value = 7144.4287
$ws.Range("K110").Value = 7144.4287
$ws.Range("M110").Value = -5099.4287

# Row 116 (G=27713)
$ws.Range("H116").Value = 2532.3157
$ws.Range("I116").Value = 2579.7334
$ws.Range("J116").Value = 2354.5
$ws.Range("K116").Value = 2579.7334
$ws.Range("L116").Value = 2354.5
$ws.Range("M116").Value = -285.7334000000001
$ws.Range("N116").Value = -6942.5

# Row 132 (G=43997)
$ws.Range("H132").Value = 8031.0884
$ws.Range("I132").Value = 7876.3438
$ws.Range("K132").Value = 23629.0314
$ws.Range("M132").Value = -21099.0314

# Row 136 (G=43999)
$ws.Range("H136").Value = 2427
$ws.Range("I136").Value = 2427
$ws.Range("K136").Value = 7281
$ws.Range("M136").Value = -4731

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (G=27713)
$ws.Range("H3").Value = 2532.3157
$ws.Range("I3").Value = 2579.7334
$ws.Range("J3").Value = 2354.5
$ws.Range("K3").Value = 2579.7334
$ws.Range("L3").Value = 2354.5
$ws.Range("M3").Value = -2465.7334
$ws.Range("N3").Value = -2582.5

# Row 94 (G=19939)
$ws.Range("H94").Value = 990.46155
$ws.Range("I94").Value = 994.08
$ws.Range("J94").Value = 900
$ws.Range("K94").Value = 994.08
$ws.Range("L94").Value = 900
$ws.Range("M94").Value = -543.08
$ws.Range("N94").Value = -1802

$ws = $wb.Worksheets.Item("CRP")
# Row 2 (G=1820)
$ws.Range("H2").Value = 251097.88
$ws.Range("I2").Value = 1254.7142
$ws.Range("K2").Value = 1254.7142
$ws.Range("M2").Value = -1141.7142

# Row 31 (G=44023)
$ws.Range("H31").Value = 3417.2856
$ws.Range("I31").Value = 2341.4
$ws.Range("K31").Value = 2341.4
$ws.Range("M31").Value = -2046.4

# Row 34 (G=44023)
$ws.Range("H34").Value = 3417.2856
$ws.Range("I34").Value = 2341.4
$ws.Range("K34").Value = 2341.4
$ws.Range("M34").Value = -2139.4

# Row 59 (G=1942)
$ws.Range("H59").Value = 19999.857
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -22290

# Row 68 (G=10611)
$ws.Range("H68").Value = 25000
$ws.Range("J68").Value = 25000
$ws.Range("L68").Value = 25000
$ws.Range("N68").Value = -26498

# Row 71 (G=10611)
$ws.Range("H71").Value = 25000
$ws.Range("J71").Value = 25000
$ws.Range("L71").Value = 75000
$ws.Range("N71").Value = -82488

# Row 74 (G=10636)
$ws.Range("H74").Value = 47754.332
$ws.Range("J74").Value = 54826
$ws.Range("L74").Value = 54826
$ws.Range("N74").Value = -56574

# Row 77 (G=10636)
$ws.Range("H77").Value = 47754.332
$ws.Range("J77").Value = 54826
$ws.Range("L77").Value = 164478
$ws.Range("N77").Value = -173214

# Row 105 (G=19928)
$ws.Range("H105").Value = 1005.26666
$ws.Range("I105").Value = 798.2727
$ws.Range("J105").Value = 1574.5
$ws.Range("K105").Value = 798.2727
$ws.Range("L105").Value = 1574.5
$ws.Range("M105").Value = 948.7273
$ws.Range("N105").Value = -5068.5

# Row 132 (G=44019)
$ws.Range("H132").Value = 2297.7778
$ws.Range("I132").Value = 2448.7368
$ws.Range("J132").Value = 1939.25
$ws.Range("K132").Value = 7346.2104
$ws.Range("L132").Value = 5817.75
$ws.Range("M132").Value = -4816.2104
$ws.Range("N132").Value = -10877.75

# Row 134 (G=44020)
$ws.Range("H134").Value = 3184.3684
$ws.Range("I134").Value = 3139.0557
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 9417.167099999999
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -6882.167099999999
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CUL")
# Row 32 (G=4731)
$ws.Range("H32").Value = 688821.75
$ws.Range("I32").Value = 202569.2
$ws.Range("J32").Value = 1094032.1
$ws.Range("K32").Value = 607707.6000000001
$ws.Range("L32").Value = 3282096.3
$ws.Range("M32").Value = -607424.6000000001
$ws.Range("N32").Value = -3282662.3

# Row 56 (G=10146)
$ws.Range("H56").Value = 8157.2144
$ws.Range("I56").Value = 8157.2144
$ws.Range("K56").Value = 8157.2144
$ws.Range("M56").Value = -7627.2144

# Row 113 (G=27843)
$ws.Range("H113").Value = 1024.4445
$ws.Range("J113").Value = 760
$ws.Range("L113").Value = 2280
$ws.Range("N113").Value = -6620

# Row 129 (G=36054)
$ws.Range("H129").Value = 113906.555
$ws.Range("I129").Value = 287013
$ws.Range("J129").Value = 3747.9092
$ws.Range("K129").Value = 861039
$ws.Range("L129").Value = 11243.7276
$ws.Range("M129").Value = -856039
$ws.Range("N129").Value = -21243.7276

# Row 131 (G=36060)
$ws.Range("H131").Value = 55601.375
$ws.Range("J131").Value = 25167.77
$ws.Range("L131").Value = 75503.31
$ws.Range("N131").Value = -85583.31

# Row 137 (G=44088)
$ws.Range("H137").Value = 5887804
$ws.Range("I137").Value = 16668114
$ws.Range("J137").Value = 7634.727
$ws.Range("K137").Value = 50004342
$ws.Range("L137").Value = 22904.181
$ws.Range("M137").Value = -49999242
$ws.Range("N137").Value = -33104.181

# Row 138 (G=44105)
$ws.Range("H138").Value = 1850
$ws.Range("I138").Value = 1000
$ws.Range("J138").Value = 2700
$ws.Range("K138").Value = 3000
$ws.Range("L138").Value = 8100
$ws.Range("M138").Value = 2140
$ws.Range("N138").Value = -18380

$ws = $wb.Worksheets.Item("GSM")
# Row 132 (G=44008)
$ws.Range("H132").Value = 3016
$ws.Range("I132").Value = 2839.3928
$ws.Range("K132").Value = 8518.178400000001
$ws.Range("M132").Value = -5988.178400000001

$ws = $wb.Worksheets.Item("LTW")
# Row 2 (G=2631)
$ws.Range("H2").Value = 100.454544
$ws.Range("J2").Value = 100.454544
$ws.Range("L2").Value = 100.454544
$ws.Range("N2").Value = -324.454544

# Row 136 (G=44060)
$ws.Range("H136").Value = 2560.0605
$ws.Range("I136").Value = 1573.5652
$ws.Range("J136").Value = 4829
$ws.Range("K136").Value = 4720.6956
$ws.Range("L136").Value = 14487
$ws.Range("M136").Value = -2170.6956
$ws.Range("N136").Value = -19587

$ws = $wb.Worksheets.Item("WVR")
# Row 11 (G=3001)
$ws.Range("H11").Value = 9500
$ws.Range("J11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("N11").Value = -10284

# Row 81 (G=12596)
$ws.Range("H81").Value = 5579.8
$ws.Range("I81").Value = 5645.3076
$ws.Range("K81").Value = 11290.6152
$ws.Range("M81").Value = -10229.6152

# Row 84 (G=12596)
$ws.Range("H84").Value = 5579.8
$ws.Range("I84").Value = 5645.3076
$ws.Range("K84").Value = 56453.076
$ws.Range("M84").Value = -51149.076
